$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.080.26'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.072.11'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -3.21%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.03'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.24'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.96%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.068.51'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.504'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.140'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.28'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.443'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000238'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.86'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.120'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.576.88'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.308.17'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.069.19'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -3.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.42'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '449.59'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.62'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.675'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.39'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.96'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.96'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('E28').Value = '  -3.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.03'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.50'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.48'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -5.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.05'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0979'
$ws.Range('D33').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.34'
$ws.Range('D34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.981'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -4.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.74'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '50.55'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0700'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0378'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.96'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('E41').Value = '  -1.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '383.38'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -7.30%  '
$ws.Range('E43').Value = '  -4.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.707.02'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -6.39%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '124.64'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.242'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.94%  '
$ws.Range('B48').Value = 'Fetch.AI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.04'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -4.44%  '
$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.04'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -5.99%  '
$ws.Range('E50').Value = '  -1.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.20'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.39%  '
